$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append the newest batch of certificate records to the log (rows 844-852).
#
# Row 844 is a standalone entry (single certificate/course) for
# IHAB SAMIR SAAD ELSAID ZANATY, formatted like the other single/"batch
# header" rows (style used across columns A-D, date column included).
#
# Rows 845-852 are a brand-new 8-row batch (one row per standard course)
# for Mohamed Ali Rajab Ali, formatted exactly like the previous batch
# (rows 836-843), reusing the same standard course-name / date strings.
# ---------------------------------------------------------------------------

# --- formatting first (copies styles only, doesn't touch cell content) ---

# Row 844: columns A-D take the "single entry" look (row 835's A-C cells
# supply that style; row 835's own A cell is reused for D844 too, since
# that row's date cell happens to carry the same plain style there).
$ws.Range("A835:C835").Copy() | Out-Null
$ws.Range("A844:C844").PasteSpecial(-4122) | Out-Null
$ws.Range("A835").Copy() | Out-Null
$ws.Range("D844").PasteSpecial(-4122) | Out-Null
$ws.Range("E843").Copy() | Out-Null
$ws.Range("E844").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(844).RowHeight = 15.75

# Rows 845-852: same look as the prior 8-row batch.
$ws.Range("A836:E843").Copy() | Out-Null
$ws.Range("A845:E852").PasteSpecial(-4122) | Out-Null
for ($r = 845; $r -le 852; $r++) {
    $ws.Rows.Item($r).RowHeight = 15.75
}

# --- values, typed in the natural order a user would enter them ---------

# Row 844
$ws.Range("A844").Value = "DSS1843"
$ws.Range("B844").Value = "IHAB SAMIR SAAD ELSAID ZANATY"
$ws.Range("C844").Value = "First Aid"
$ws.Range("D844").Value = "20-12-2024"
$ws.Range("E844").Value = 1

# Certificate numbers for the new batch (column A, filled down first)
$ws.Range("A845").Value = "DSS1844"
$ws.Range("A846").Value = "DSS1845"
$ws.Range("A847").Value = "DSS1846"
$ws.Range("A848").Value = "DSS1847"
$ws.Range("A849").Value = "DSS1848"
$ws.Range("A850").Value = "DSS1849"
$ws.Range("A851").Value = "DSS1850"
$ws.Range("A852").Value = "DSS1851"

# Name, same for all 8 rows of the batch (column B, filled across next)
$ws.Range("B845:B852").Value = "Mohamed Ali Rajab Ali"

# Course names and dates (columns C & D), then the result flag (column E).
# The D-column cells are formatted with a date number format (style 45),
# but the dates here are kept as plain text (as in every other batch in
# this sheet), so a leading apostrophe forces text entry instead of
# letting Excel auto-convert the typed string into a real date value.
$ws.Range("C845").Value = "30 Hours Construction Safety & Health"
$ws.Range("D845").Value = "'05-12-2024"
$ws.Range("E845").Value = 1

$ws.Range("C846").Value = "30 Hours G. Industry Safety & Health"
$ws.Range("D846").Value = "'10-12-2024"
$ws.Range("E846").Value = 1

$ws.Range("C847").Value = "Electrical Safety & LOTO"
$ws.Range("D847").Value = "'06-12-2024"
$ws.Range("E847").Value = 1

$ws.Range("C848").Value = "Fire Marshal"
$ws.Range("D848").Value = "'03-12-2024"
$ws.Range("E848").Value = 1

$ws.Range("C849").Value = "Scaffold Competent Person"
$ws.Range("D849").Value = "'01-12-2024"
$ws.Range("E849").Value = 1

$ws.Range("C850").Value = "Lifting & Rigging Competent Person"
$ws.Range("D850").Value = "'02-12-2024"
$ws.Range("E850").Value = 1

$ws.Range("C851").Value = "Health & Safety Risk Assessment"
$ws.Range("D851").Value = "'07-12-2024"
$ws.Range("E851").Value = 1

$ws.Range("C852").Value = "Safety Management System & PTW"
$ws.Range("D852").Value = "'08-12-2024"
$ws.Range("E852").Value = 1

# Reflect the user's final selection/scroll position after typing the rows.
$ws.Range("C853").Select()
$excel.ActiveWindow.ScrollRow = 835
$excel.ActiveWindow.ScrollColumn = 1
